# Adactin App changes 02-11-2022
#
# On the "Adactin" sheet:
#  - G2 used to hold the hyperlinked text "Adactin@123"; it is replaced by
#    the plain number 123456 (and its hyperlink is removed).
#  - H2 used to read "NAVEENREDMINOTE10"; it is replaced by "NaveenSelenium".
#  - The sheet's active selection moves to G17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Adactin")

# Drop the hyperlink anchored on G2 before it becomes a plain numeric value.
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$G$2') {
        $hl.Delete()
        break
    }
}

# G2: "Adactin@123" (hyperlinked text) -> 123456 (plain number)
$ws.Range("G2").Value = 123456

# H2: "NAVEENREDMINOTE10" -> "NaveenSelenium"
$ws.Range("H2").Value = "NaveenSelenium"

# Move the active selection on the Adactin sheet to G17
[void]$ws.Activate()
[void]$ws.Range("G17").Select()
